# Auto-generated Excel COM-interop script
# Updates market-price / profit figures on several sheets (ALC, BSM, CRP, CUL, GSM, LTW, WVR)
# as pulled by the scheduled Sheets runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 12353535
$ws.Range("I40").Value = 5233
$ws.Range("K40").Value = 5233
$ws.Range("M40").Value = -5058
$ws.Range("H116").Value = 3832.6875
$ws.Range("I116").Value = 3773.0715
$ws.Range("K116").Value = 3773.0715
$ws.Range("M116").Value = -331.0715
$ws.Range("H138").Value = 5347.222
$ws.Range("I138").Value = 3753.8125
$ws.Range("J138").Value = 6018.1313
$ws.Range("K138").Value = 11261.4375
$ws.Range("L138").Value = 18054.3939
$ws.Range("M138").Value = -6121.4375
$ws.Range("N138").Value = -28334.3939

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14311551
$ws.Range("I134").Value = 14311551
$ws.Range("K134").Value = 42934653
$ws.Range("M134").Value = -42932118

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16678827
$ws.Range("I58").Value = 25017482
$ws.Range("K58").Value = 25017482
$ws.Range("M58").Value = -25017279
$ws.Range("H86").Value = 7099.6665
$ws.Range("I86").Value = 5549.625
$ws.Range("K86").Value = 5549.625
$ws.Range("M86").Value = -4426.625
$ws.Range("H89").Value = 7099.6665
$ws.Range("I89").Value = 5549.625
$ws.Range("K89").Value = 27748.125
$ws.Range("M89").Value = -22132.125
$ws.Range("H97").Value = 29999.5
$ws.Range("J97").Value = 29999.5
$ws.Range("L97").Value = 29999.5
$ws.Range("N97").Value = -31981.5
$ws.Range("H99").Value = 3496.6667
$ws.Range("I99").Value = 3495
$ws.Range("J99").Value = 3497.5
$ws.Range("K99").Value = 3495
$ws.Range("L99").Value = 3497.5
$ws.Range("M99").Value = -1997
$ws.Range("N99").Value = -6493.5
$ws.Range("H126").Value = 3496.6667
$ws.Range("I126").Value = 3495
$ws.Range("J126").Value = 3497.5
$ws.Range("K126").Value = 10485
$ws.Range("L126").Value = 10492.5
$ws.Range("M126").Value = -8015
$ws.Range("N126").Value = -15432.5
$ws.Range("H132").Value = 58825030
$ws.Range("I132").Value = 58825030
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 176475090
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -176472560
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 14707845
$ws.Range("I134").Value = 25001850
$ws.Range("J134").Value = 2125.2856
$ws.Range("K134").Value = 75005550
$ws.Range("L134").Value = 6375.8568
$ws.Range("M134").Value = -75003015
$ws.Range("N134").Value = -11445.8568
$ws.Range("H136").Value = 16678827
$ws.Range("I136").Value = 25017482
$ws.Range("K136").Value = 75052446
$ws.Range("M136").Value = -75049896
$ws.Range("H140").Value = 57713.57
$ws.Range("J140").Value = 57713.57
$ws.Range("L140").Value = 57713.57
$ws.Range("N140").Value = -68073.57000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 520164.44
$ws.Range("I4").Value = 700947.0600000001
$ws.Range("K4").Value = 2102841.18
$ws.Range("M4").Value = -2102729.18
$ws.Range("H14").Value = 476
$ws.Range("I14").Value = 476
$ws.Range("K14").Value = 1428
$ws.Range("M14").Value = -1255
$ws.Range("H23").Value = 803.4737
$ws.Range("J23").Value = 1265
$ws.Range("L23").Value = 3795
$ws.Range("N23").Value = -4265
$ws.Range("H34").Value = 1074.6
$ws.Range("I34").Value = 361.75
$ws.Range("K34").Value = 1085.25
$ws.Range("M34").Value = -1001.25
$ws.Range("H39").Value = 2548.3333
$ws.Range("I39").Value = 1625
$ws.Range("J39").Value = 7165
$ws.Range("K39").Value = 4875
$ws.Range("L39").Value = 21495
$ws.Range("M39").Value = -4581
$ws.Range("N39").Value = -22083
$ws.Range("H47").Value = 626.2857
$ws.Range("I47").Value = 724
$ws.Range("K47").Value = 2172
$ws.Range("M47").Value = -1741
$ws.Range("H55").Value = 1962.25
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1962.25
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 5886.75
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -6240.75
$ws.Range("H68").Value = 6474.5
$ws.Range("I68").Value = 6400
$ws.Range("J68").Value = 6499.3335
$ws.Range("K68").Value = 19200
$ws.Range("L68").Value = 19498.0005
$ws.Range("M68").Value = -18389
$ws.Range("N68").Value = -21120.0005
$ws.Range("H71").Value = 6474.5
$ws.Range("I71").Value = 6400
$ws.Range("J71").Value = 6499.3335
$ws.Range("K71").Value = 57600
$ws.Range("L71").Value = 58494.0015
$ws.Range("M71").Value = -53544
$ws.Range("N71").Value = -66606.0015

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5532.6665
$ws.Range("I70").Value = 5532.6665
$ws.Range("K70").Value = 5532.6665
$ws.Range("M70").Value = -5262.6665
$ws.Range("H73").Value = 5532.6665
$ws.Range("I73").Value = 5532.6665
$ws.Range("K73").Value = 5532.6665
$ws.Range("M73").Value = -4596.6665
$ws.Range("H75").Value = 150000
$ws.Range("J75").Value = 150000
$ws.Range("L75").Value = 150000
$ws.Range("N75").Value = -151748
$ws.Range("H78").Value = 150000
$ws.Range("J78").Value = 150000
$ws.Range("L78").Value = 450000
$ws.Range("N78").Value = -458736
$ws.Range("H80").Value = 4319
$ws.Range("I80").Value = 4319
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4319
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3321
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 4319
$ws.Range("I83").Value = 4319
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 21595
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -16603
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 5016109.5
$ws.Range("I132").Value = 5697214.5
$ws.Range("J132").Value = 21338
$ws.Range("K132").Value = 17091643.5
$ws.Range("L132").Value = 64014
$ws.Range("M132").Value = -17089113.5
$ws.Range("N132").Value = -69074

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 456.53333
$ws.Range("J55").Value = 683.5
$ws.Range("L55").Value = 683.5
$ws.Range("N55").Value = -1029.5
$ws.Range("H61").Value = 5067.25
$ws.Range("I61").Value = 5186.0557
$ws.Range("K61").Value = 5186.0557
$ws.Range("M61").Value = -4984.0557
$ws.Range("H113").Value = 5067.25
$ws.Range("I113").Value = 5186.0557
$ws.Range("K113").Value = 5186.0557
$ws.Range("M113").Value = -3016.0557
$ws.Range("H122").Value = 3419
$ws.Range("I122").Value = 3419
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10257
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7807
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 50024000
$ws.Range("I132").Value = 62528750
$ws.Range("K132").Value = 187586250
$ws.Range("M132").Value = -187583720
$ws.Range("H136").Value = 2052.6191
$ws.Range("J136").Value = 2635.2856
$ws.Range("L136").Value = 7905.8568
$ws.Range("N136").Value = -13005.8568

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7399.2
$ws.Range("J62").Value = 8499
$ws.Range("L62").Value = 8499
$ws.Range("N62").Value = -9747
$ws.Range("H65").Value = 7399.2
$ws.Range("J65").Value = 8499
$ws.Range("L65").Value = 42495
$ws.Range("N65").Value = -48735
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 3209.3333
$ws.Range("J100").Value = 1000
$ws.Range("L100").Value = 2000
$ws.Range("N100").Value = -3082
$ws.Range("H132").Value = 9437554
$ws.Range("I132").Value = 11629393
$ws.Range("K132").Value = 34888179
$ws.Range("M132").Value = -34885649
